$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BUMP_DOWN")

# ---------------------------------------------------------------------
# Row 8: column letter labels (B8:AA8), mirrored left-to-right relative
# to the "Full" sheet columns, plus an empty styled corner cell (A8).
# ---------------------------------------------------------------------
$colLabels = @("AF","AE","AD","AC","AB","AA","Y","W","V","U","T","R","P","N","M","L","K","J","H","G","F","E","D","C","B","A")

$ws.Range("A8").Value = ""

$startCol = 2 # column B
for ($i = 0; $i -lt $colLabels.Length; $i++) {
    $ws.Cells.Item(8, $startCol + $i).Value = $colLabels[$i]
}

# ---------------------------------------------------------------------
# Column A (rows 9:34): row numbers 1..26
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 26; $i++) {
    $ws.Cells.Item(9 + $i, 1).Value = $i + 1
}

# ---------------------------------------------------------------------
# Formatting - corner cell A8 uses a smaller Comic Sans MS (12pt)
# ---------------------------------------------------------------------
$corner = $ws.Range("A8")
$corner.Font.Name = "Comic Sans MS"
$corner.Font.Size = 12
$corner.HorizontalAlignment = -4108
$corner.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Formatting - the row/column labels use a larger Comic Sans MS (16pt)
# (applied per-area; multi-area Range font assignment only affects the
# first area in this runtime, so the two areas are handled separately)
# ---------------------------------------------------------------------
$colHeader = $ws.Range("B8:AA8")
$colHeader.Font.Name = "Comic Sans MS"
$colHeader.Font.Size = 16
$colHeader.HorizontalAlignment = -4108
$colHeader.VerticalAlignment = -4108

$rowHeader = $ws.Range("A9:A34")
$rowHeader.Font.Name = "Comic Sans MS"
$rowHeader.Font.Size = 16
$rowHeader.HorizontalAlignment = -4108
$rowHeader.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# View state - user zoomed in and moved the selection while editing
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("H6").Select() | Out-Null
